$d = $word.ActiveDocument
$d.Content.Find.Execute("Documentação  do Gerador de", $true, $false, $false, $false, $false, $true, 1, $false, "Documentação do Gerador de", 2)
